$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency table with refreshed prices / volume figures.
# Cells whose new text looks like a plain number (e.g. "568.30", "0.420")
# are first switched to Text format so Excel keeps the exact original
# string (trailing zeros, etc.) instead of coercing it into a float.

$ws.Range("D2").Value = "63.728.54"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "3.406.44"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.30"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.39"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.412.71"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  -9.18%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.420"
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("D13").Value = "4.001.42"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.94"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -8.90%  "
$ws.Range("D17").Value = "63.863.66"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "3.400.06"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  -4.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.00"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.80"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.42"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.513"
$ws.Range("E25").Value = "  -6.52%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  -5.44%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("E31").Value = "  -6.08%  "
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.74"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.29"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.832"
$ws.Range("E37").Value = "  +7.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  -6.47%  "
$ws.Range("D39").Value = "2.809.15"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0727"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.94"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.49"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.59"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.97"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0301"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.54"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +5.11%  "
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  -4.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.27"
$ws.Range("E51").Value = "  -4.16%  "
